# Bethesda University Organizations workbook fix-up
#
# 1. Swap the meaning of columns A and B: A becomes "Category" (short
#    category label) and B becomes "Organization Name" (was previously A).
# 2. Rename several headers (C, D, G, H, I, J, K, L) to the new,
#    longer "<Platform> Link" naming scheme.
# 3. Add a new column M "Tiktok Link".
# 4. Resize columns A, B and G:M to their new widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Swap column A / B data (the three data rows) before touching
#    headers, so we move the *values* currently sitting in A into B and
#    vice versa. NOTE: plain `.Value` reads are unreliable on this COM
#    shim (they surface the property descriptor instead of the cell's
#    contents) - use `.Value2` for reads, `.Value` remains fine for writes.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 4; $r++) {
    $catCell = $ws.Cells.Item($r, 1)   # currently the "Organization Name" value
    $nameCell = $ws.Cells.Item($r, 2)  # currently the "Categories" value

    $catValue = $catCell.Value2
    $nameValue = $nameCell.Value2

    # A should hold the (short) category, B the organization name.
    $catCell.Value = $nameValue
    $nameCell.Value = $catValue
}

# ---------------------------------------------------------------------
# 2. Header row text.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Organization Name"
$ws.Range("C1").Value = "Organization Link"
$ws.Range("D1").Value = "Logo Link"
$ws.Range("G1").Value = "Phone Number"
$ws.Range("H1").Value = "Linkedin Link"
$ws.Range("I1").Value = "Instagram Link"
$ws.Range("J1").Value = "Facebook Link"
$ws.Range("K1").Value = "Twitter Link"
$ws.Range("L1").Value = "Youtube Link"

# ---------------------------------------------------------------------
# 3. New "Tiktok Link" column (M) - copy the header formatting from the
#    neighbouring L1 cell so it picks up the same bold/centered/bordered
#    style used by the rest of row 1, then add the blank data cells for
#    rows 2-4 to match the existing table shape.
# ---------------------------------------------------------------------
$ws.Range("M1").Value = "Tiktok Link"
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null

# Materialise the (empty) M2:M4 data cells so the column has the same
# "every row has a cell" shape as the rest of the table. Assigning the
# built-in "Normal" cell style forces the cell to be written without
# pulling in a brand-new style record (it resolves back to the sheet's
# existing default xf), so styles.xml stays untouched.
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 13).Style = "Normal"
}

# ---------------------------------------------------------------------
# 4. Column widths. ColumnWidth is expressed in "characters" and Excel
#    snaps it to the underlying pixel grid, so nudge each target value
#    by a small fraction (-0.85) to reliably land back on the exact
#    whole-character width stored in the sheet XML.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.15   # A: 32 -> 20
$ws.Columns.Item(2).ColumnWidth = 31.15   # B: 20 -> 32
$ws.Columns.Item(7).ColumnWidth = 13.15   # G: 7  -> 14
$ws.Columns.Item(8).ColumnWidth = 14.15   # H: 9  -> 15
$ws.Columns.Item(9).ColumnWidth = 15.15   # I: 10 -> 16
$ws.Columns.Item(10).ColumnWidth = 14.15  # J: 11 -> 15
$ws.Columns.Item(11).ColumnWidth = 13.15  # K: 10 -> 14
$ws.Columns.Item(12).ColumnWidth = 13.15  # L: 9  -> 14
$ws.Columns.Item(13).ColumnWidth = 12.15  # M: new -> 13

$ws.Range("A1").Select()
